$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 60 (pushes the existing rows 60-95 down to 61-96)
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(60, 1).Value  = 11
$ws.Cells.Item(60, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(60, 3).Value  = "Bíobío"
$ws.Cells.Item(60, 4).Value  = 45126
$ws.Cells.Item(60, 5).Value  = 8
$ws.Cells.Item(60, 6).Value  = 100112013
$ws.Cells.Item(60, 7).Value  = "Alcachofa"
$ws.Cells.Item(60, 8).Value  = "Argentina(o)"
$ws.Cells.Item(60, 9).Value  = "Primera"
$ws.Cells.Item(60, 10).Value = 110
$ws.Cells.Item(60, 11).Value = 14000
$ws.Cells.Item(60, 12).Value = 15000
$ws.Cells.Item(60, 13).Value = 14545
$ws.Cells.Item(60, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(60, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(60, 16).Value = 291
$ws.Cells.Item(60, 17).Value = 50
$ws.Cells.Item(60, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(61, 4).NumberFormat
